# Optuna Attempt (go back with original)
# Updates forecast values on "Forecast Comparison" sheet and summary totals
# on the "Summary" sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet updates ---

# Row 2
$ws1.Range("D2").Value = 18
$ws1.Range("H2").Value = 0.38
$ws1.Range("I2").Value = "High"
$ws1.Range("L2").Value = 1.11

# Row 3
$ws1.Range("D3").Value = 22
$ws1.Range("L3").Value = 0.89

# Row 4
$ws1.Range("D4").Value = 24
$ws1.Range("L4").Value = 0.97

# Row 5
$ws1.Range("D5").Value = 31
$ws1.Range("L5").Value = 1.05

# Row 6
$ws1.Range("D6").Value = 29
$ws1.Range("L6").Value = 1.11

# Row 7
$ws1.Range("D7").Value = 14
$ws1.Range("L7").Value = 0.93

# Row 8
$ws1.Range("L8").Value = 1.14

# Row 9
$ws1.Range("L9").Value = 1.06

# Row 10
$ws1.Range("D10").Value = 14
$ws1.Range("L10").Value = 1.07

# Row 11
$ws1.Range("D11").Value = 14
$ws1.Range("L11").Value = 0.89

# Row 12
$ws1.Range("L12").Value = 1.04

# Row 13
$ws1.Range("D13").Value = 14
$ws1.Range("L13").Value = 0.9

# Row 14
$ws1.Range("D14").Value = 14
$ws1.Range("L14").Value = 0.95

# Row 15
$ws1.Range("D15").Value = 14
$ws1.Range("L15").Value = 0.95

# Row 16
$ws1.Range("D16").Value = 14
$ws1.Range("L16").Value = 1.11

# Row 17
$ws1.Range("D17").Value = 14
$ws1.Range("L17").Value = 0.83

# --- Summary sheet updates ---
# These cells hold their numbers as text, so force text formatting before
# assigning the value to avoid Excel auto-converting them to numerics.

$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "286"

$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "171"

$ws2.Range("B11").NumberFormat = "@"
$ws2.Range("B11").Value = "96"

$ws2.Range("B12").NumberFormat = "@"
$ws2.Range("B12").Value = "31"
